# Auto-generated Excel COM-interop script applying the cryptos-list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '29.460.23'
Set-TextCell $ws 'E2' '  -2.20%  '

Set-TextCell $ws 'D3' '1.989.77'
Set-TextCell $ws 'E3' '  -6.02%  '

Set-TextCell $ws 'E4' '  -0.05%  '

Set-TextCell $ws 'D5' '329.94'
Set-TextCell $ws 'E5' '  -4.75%  '

Set-TextCell $ws 'E6' '  +0.03%  '

Set-TextCell $ws 'D7' '0.4979'
Set-TextCell $ws 'E7' '  -4.45%  '

Set-TextCell $ws 'D8' '0.4211'
Set-TextCell $ws 'E8' '  -5.78%  '

Set-TextCell $ws 'D9' '52.00'
Set-TextCell $ws 'E9' '  -4.15%  '

Set-TextCell $ws 'D10' '0.08889'
Set-TextCell $ws 'E10' '  -5.07%  '

Set-TextCell $ws 'D11' '1.120'
Set-TextCell $ws 'E11' '  -5.49%  '

Set-TextCell $ws 'D12' '23.33'
Set-TextCell $ws 'E12' '  -7.91%  '

Set-TextCell $ws 'D13' '8.059'
Set-TextCell $ws 'E13' '  -7.18%  '

Set-TextCell $ws 'D14' '1.968.33'
Set-TextCell $ws 'E14' '  -7.29%  '

Set-TextCell $ws 'D15' '6.499'
Set-TextCell $ws 'E15' '  -6.85%  '

Set-TextCell $ws 'D16' '96.11'
Set-TextCell $ws 'E16' '  -6.34%  '

Set-TextCell $ws 'E17' '  -0.04%  '

Set-TextCell $ws 'E18' '  -5.66%  '

Set-TextCell $ws 'D19' '0.06617'
Set-TextCell $ws 'E19' '  -1.26%  '

Set-TextCell $ws 'D20' '19.71'
Set-TextCell $ws 'E20' '  -8.65%  '

Set-TextCell $ws 'E21' '  -0.08%  '

Set-TextCell $ws 'D22' '5.957'
Set-TextCell $ws 'E22' '  -5.59%  '

Set-TextCell $ws 'D23' '29.469.64'
Set-TextCell $ws 'E23' '  -2.24%  '

Set-TextCell $ws 'D24' '11.85'
Set-TextCell $ws 'E24' '  -7.12%  '

Set-TextCell $ws 'D25' '2.273'
Set-TextCell $ws 'E25' '  -2.49%  '

Set-TextCell $ws 'D26' '157.41'
Set-TextCell $ws 'E26' '  -3.51%  '

Set-TextCell $ws 'D27' '20.57'
Set-TextCell $ws 'E27' '  -7.26%  '

Set-TextCell $ws 'D28' '6.538'
Set-TextCell $ws 'E28' '  -4.27%  '

Set-TextCell $ws 'D29' '2.334'
Set-TextCell $ws 'E29' '  -8.63%  '

Set-TextCell $ws 'D30' '127.88'
Set-TextCell $ws 'E30' '  -4.77%  '

Set-TextCell $ws 'D31' '1.052'
Set-TextCell $ws 'E31' '  -9.16%  '

Set-TextCell $ws 'D32' '0.09937'
Set-TextCell $ws 'E32' '  -6.16%  '

Set-TextCell $ws 'D33' '1.568'
Set-TextCell $ws 'E33' '  -12.33%  '

Set-TextCell $ws 'D34' '5.838'
Set-TextCell $ws 'E34' '  -7.21%  '

Set-TextCell $ws 'D35' '3.792'
Set-TextCell $ws 'E35' '  -4.43%  '

Set-TextCell $ws 'D36' '9.604'
Set-TextCell $ws 'E36' '  -10.82%  '

Set-TextCell $ws 'D37' '0.02451'
Set-TextCell $ws 'E37' '  -7.38%  '

Set-TextCell $ws 'D38' '0.06346'

Set-TextCell $ws 'E39' '  -3.64%  '

Set-TextCell $ws 'B40' 'TheSandbox'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws 'D40' '0.6513'
Set-TextCell $ws 'E40' '  -8.73%  '

Set-TextCell $ws 'B41' 'Aptos'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D41' '11.75'
Set-TextCell $ws 'E41' '  -8.02%  '

Set-TextCell $ws 'D42' '0.2066'
Set-TextCell $ws 'E42' '  -8.20%  '

Set-TextCell $ws 'E43' '  -0.03%  '

Set-TextCell $ws 'D44' '0.6349'
Set-TextCell $ws 'E44' '  -8.94%  '

Set-TextCell $ws 'D45' '2.219'
Set-TextCell $ws 'E45' '  -7.58%  '

Set-TextCell $ws 'D46' '13.37'
Set-TextCell $ws 'E46' '  -9.05%  '

Set-TextCell $ws 'D47' '1.269'
Set-TextCell $ws 'E47' '  +0.30%  '

Set-TextCell $ws 'D48' '3.531'
Set-TextCell $ws 'E48' '  -2.57%  '

Set-TextCell $ws 'D49' '0.00000000331'
Set-TextCell $ws 'E49' '  -5.44%  '

Set-TextCell $ws 'D50' '0.06982'
Set-TextCell $ws 'E50' '  -3.16%  '

Set-TextCell $ws 'E51' '  -6.65%  '

Write-Output "Applied 95 cell updates."
